$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the Price/Volume columns so values like "1.007" are not
# auto-converted to numbers by Excel (matches original inline-string content).
$ws.Range("D2:E51").NumberFormat = "@"

$ws.Range('D2').Value = '26.202.64'
$ws.Range('E2').Value = '  -4.40%  '
$ws.Range('D3').Value = '1.659.26'
$ws.Range('E3').Value = '  -3.05%  '
$ws.Range('D4').Value = '1.007'
$ws.Range('E4').Value = '  +0.24%  '
$ws.Range('D5').Value = '218.16'
$ws.Range('E5').Value = '  -2.75%  '
$ws.Range('D6').Value = '0.5154'
$ws.Range('E6').Value = '  -3.51%  '
$ws.Range('D7').Value = '1.008'
$ws.Range('E7').Value = '  +0.28%  '
$ws.Range('D8').Value = '0.2567'
$ws.Range('E8').Value = '  -3.96%  '
$ws.Range('E9').Value = '  -3.33%  '
$ws.Range('D10').Value = '19.84'
$ws.Range('E10').Value = '  -5.21%  '
$ws.Range('D11').Value = '0.07795'
$ws.Range('E11').Value = '  +2.04%  '
$ws.Range('D12').Value = '1.663.96'
$ws.Range('E12').Value = '  -3.52%  '
$ws.Range('D13').Value = '4.302'
$ws.Range('E13').Value = '  -5.60%  '
$ws.Range('D14').Value = '1.885.94'
$ws.Range('E14').Value = '  -3.17%  '
$ws.Range('D15').Value = '0.5532'
$ws.Range('E15').Value = '  -4.06%  '
$ws.Range('D16').Value = '0.0₅8028'
$ws.Range('E16').Value = '  -1.71%  '
$ws.Range('D17').Value = '64.21'
$ws.Range('E17').Value = '  -5.42%  '
$ws.Range('D18').Value = '26.221.26'
$ws.Range('E18').Value = '  -4.31%  '
$ws.Range('E19').Value = '  +0.33%  '
$ws.Range('D20').Value = '209.96'
$ws.Range('E20').Value = '  -3.16%  '
$ws.Range('D21').Value = '4.394'
$ws.Range('E21').Value = '  -5.75%  '
$ws.Range('D22').Value = '10.07'
$ws.Range('E22').Value = '  -3.71%  '
$ws.Range('D23').Value = '5.880'
$ws.Range('E23').Value = '  -1.59%  '
$ws.Range('D24').Value = '1.007'
$ws.Range('E24').Value = '  +0.24%  '
$ws.Range('D25').Value = '143.47'
$ws.Range('E25').Value = '  +0.62%  '
$ws.Range('D26').Value = '1.765'
$ws.Range('E26').Value = '  +1.97%  '
$ws.Range('E27').Value = '  -4.48%  '
$ws.Range('D28').Value = '6.966'
$ws.Range('E28').Value = '  -4.38%  '
$ws.Range('D29').Value = '15.73'
$ws.Range('E29').Value = '  -3.35%  '
$ws.Range('D30').Value = '0.05249'
$ws.Range('E30').Value = '  -2.88%  '
$ws.Range('E31').Value = '  -2.66%  '
$ws.Range('E32').Value = '  -3.47%  '
$ws.Range('D33').Value = '3.209'
$ws.Range('E33').Value = '  -6.34%  '
$ws.Range('D34').Value = '1.564'
$ws.Range('E34').Value = '  -4.74%  '
$ws.Range('D35').Value = '2.752'
$ws.Range('E35').Value = '  -4.43%  '
$ws.Range('D36').Value = '2.363'
$ws.Range('E36').Value = '  -2.15%  '
$ws.Range('E37').Value = '  -2.66%  '
$ws.Range('D38').Value = '0.5709'
$ws.Range('E38').Value = '  -2.32%  '
$ws.Range('D39').Value = '1.157.75'
$ws.Range('E39').Value = '  +10.87%  '
$ws.Range('E40').Value = '  -2.72%  '
$ws.Range('E41').Value = '  +0.29%  '
$ws.Range('D42').Value = '0.8394'
$ws.Range('E42').Value = '  -0.04%  '
$ws.Range('D43').Value = '5.667'
$ws.Range('E43').Value = '  -3.32%  '
$ws.Range('D44').Value = '99.94'
$ws.Range('E44').Value = '  -0.88%  '
$ws.Range('D45').Value = '1.795.96'
$ws.Range('E45').Value = '  -3.18%  '
$ws.Range('D46').Value = '0.0₈111'
$ws.Range('E46').Value = '  -6.91%  '
$ws.Range('D47').Value = '0.4507'
$ws.Range('E47').Value = '  -0.27%  '
$ws.Range('D48').Value = '56.00'
$ws.Range('E48').Value = '  -3.40%  '
$ws.Range('E49').Value = '  +0.59%  '
$ws.Range('D50').Value = '7.906'
$ws.Range('E50').Value = '  -2.40%  '
$ws.Range('D51').Value = '0.05093'
$ws.Range('E51').Value = '  -2.88%  '

# Remove the temporary formatting so the cell style matches the original (no explicit style).
$ws.Range("D2:E51").ClearFormats()
